$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "donneesActuelles" column (U, 21st),
# shifting donneesActuelles -> V and anomalies -> W.
$ws.Columns.Item(21).Insert()

# New header for the inserted column.
$ws.Range("U1").Value = "source"

# Fill the new "source" column for the 4 data rows.
$ws.Range("U2").Value = "data.gouv.fr_aife"
$ws.Range("U3").Value = "data.gouv.fr_aife"
$ws.Range("U4").Value = "data.gouv.fr_aife"
$ws.Range("U5").Value = "data.gouv.fr_aife"

# "donneesActuelles" (now column V) switches from numeric 0/1 to textual non/oui.
$ws.Range("V2").Value = "non"
$ws.Range("V3").Value = "non"
$ws.Range("V4").Value = "oui"
$ws.Range("V5").Value = "oui"

# Column width for the newly inserted "source" column (U). The previously
# existing columns (now V/W) keep their original widths automatically since
# the column insert shifts them along with their <col> definitions.
# (This runtime snaps ColumnWidth to 1/6-character increments like real
# Excel; 14 is the input that lands closest to the target 14.9 width.)
$ws.Columns.Item(21).ColumnWidth = 14
